$d = $word.ActiveDocument

# Locate the sentence that gets the new trailing "(Changed main)" annotation.
$hit = $d.Content
$hit.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)

# Range covering just the matched text (no paragraph mark), which we will
# replace in-place via raw OOXML so the three appended runs come out as
# separate <w:r> elements (matching the target diff) instead of being
# coalesced into one run the way plain InsertAfter calls would be.
$target = $d.Range($hit.Start, $hit.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
       '<pkg:xmlData>' + `
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
       '<w:body><w:p>' + `
       '<w:r><w:t>This is a Microsoft word document.</w:t></w:r>' + `
       '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' + `
       '<w:r><w:t>Changed main</w:t></w:r>' + `
       '<w:r><w:t>)</w:t></w:r>' + `
       '</w:p></w:body></w:document>' + `
       '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)
